# "Rest Assured codes added" — update the Control Flag on the HomePage sheet
# from "No" to "Yes" for the Coorg row, and move the cell cursor to A4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HomePage")
$ws.Activate()

$ws.Range("A3").Value = "Yes"
$ws.Range("A4").Select()
